$d = $word.ActiveDocument

# Update the date heading paragraph.
$d.Paragraphs.Item(1).Range.Text = "2026-02-25 Wednesday"

# Update the multiplication table entries.
# The table has 20 rows x 5 cols but only rows 1, 5, 10, 15 and 20
# contain text. We address each cell explicitly by row/col and set its
# Range.Text directly (rather than using Find.Execute scoped to the cell)
# because some source values are duplicated across cells in the same row
# (e.g. row 5 col 1 and row 5 col 2 both start as "93x83=7719"), and a
# text search/replace is not reliably confined to the target cell's range
# when an identical run of text also exists elsewhere in the document.
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; New="79×26=2054"},
    @{Row=1;  Col=2; New="43×92=3956"},
    @{Row=1;  Col=3; New="47×90=4230"},
    @{Row=1;  Col=4; New="21×30=630"},
    @{Row=1;  Col=5; New="62×13=806"},

    @{Row=5;  Col=1; New="61×54=3294"},
    @{Row=5;  Col=2; New="13×36=468"},
    @{Row=5;  Col=3; New="54×76=4104"},
    @{Row=5;  Col=4; New="68×41=2788"},
    @{Row=5;  Col=5; New="34×62=2108"},

    @{Row=10; Col=1; New="15×66=990"},
    @{Row=10; Col=2; New="41×53=2173"},
    @{Row=10; Col=3; New="85×42=3570"},
    @{Row=10; Col=4; New="71×74=5254"},
    @{Row=10; Col=5; New="38×68=2584"},

    @{Row=15; Col=1; New="95×52=4940"},
    @{Row=15; Col=2; New="53×90=4770"},
    @{Row=15; Col=3; New="91×74=6734"},
    @{Row=15; Col=4; New="40×88=3520"},
    @{Row=15; Col=5; New="38×49=1862"},

    @{Row=20; Col=1; New="45×78=3510"},
    @{Row=20; Col=2; New="75×80=6000"},
    @{Row=20; Col=3; New="48×28=1344"},
    @{Row=20; Col=4; New="45×18=810"},
    @{Row=20; Col=5; New="72×81=5832"}
)

foreach ($u in $updates) {
    $t.Cell($u.Row, $u.Col).Range.Text = $u.New
}
